$wb = $excel.ActiveWorkbook

# The "About" sheet is the first sheet in the workbook
$ws = $wb.Worksheets.Item("About")

# Add a date stamp in cell C1 (Excel serial date 44307 = 2021-04-21)
$cell = $ws.Range("C1")
$cell.Value = 44307
$cell.NumberFormat = "mm-dd-yy"
